$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2058823529411765
$ws.Range("C2").Value = 0.5203619909502263
$ws.Range("J2").Value = 0.01809954751131222
$ws.Range("P2").Value = 0.1380090497737557
$ws.Range("S2").Value = 0.1176470588235294
$ws.Range("B3").Value = 0.008583690987124463
$ws.Range("C3").Value = 0.02145922746781116
$ws.Range("J3").Value = 0.03862660944206009
$ws.Range("O3").Value = 0.004291845493562232
$ws.Range("P3").Value = 0.7253218884120172
$ws.Range("S3").Value = 0.2017167381974249
$ws.Range("J4").Value = 0.08928571428571429
$ws.Range("O4").Value = 0.01785714285714286
$ws.Range("P4").Value = 0.7321428571428571
$ws.Range("S4").Value = 0.1607142857142857
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("E6").Value = 0.004166666666666667
$ws.Range("F6").Value = 0.07916666666666666
$ws.Range("J6").Value = 0.1958333333333333
$ws.Range("O6").Value = 0.04166666666666666
$ws.Range("Q6").Value = 0.1875
$ws.Range("R6").Value = 0.075
$ws.Range("S6").Value = 0.3166666666666667
$ws.Range("B7").Value = 0.1226415094339623
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.07547169811320754
$ws.Range("J7").Value = 0.1367924528301887
$ws.Range("O7").Value = 0.0660377358490566
$ws.Range("Q7").Value = 0.1226415094339623
$ws.Range("R7").Value = 0.05188679245283019
$ws.Range("S7").Value = 0.4056603773584906
$ws.Range("B8").Value = 0.1605691056910569
$ws.Range("D8").Value = 0.02845528455284553
$ws.Range("F8").Value = 0.05691056910569105
$ws.Range("J8").Value = 0.1097560975609756
$ws.Range("O8").Value = 0.0508130081300813
$ws.Range("Q8").Value = 0.1951219512195122
$ws.Range("R8").Value = 0.06707317073170732
$ws.Range("S8").Value = 0.3313008130081301
$ws.Range("B9").Value = 0.136
$ws.Range("D9").Value = 0.016
$ws.Range("F9").Value = 0.07199999999999999
$ws.Range("J9").Value = 0.104
$ws.Range("O9").Value = 0.048
$ws.Range("Q9").Value = 0.208
$ws.Range("R9").Value = 0.048
$ws.Range("S9").Value = 0.368
$ws.Range("B10").Value = 0.1396363636363636
$ws.Range("D10").Value = 0.02327272727272727
$ws.Range("F10").Value = 0.06254545454545454
$ws.Range("J10").Value = 0.1163636363636364
$ws.Range("O10").Value = 0.03854545454545454
$ws.Range("Q10").Value = 0.2210909090909091
$ws.Range("R10").Value = 0.06909090909090909
$ws.Range("S10").Value = 0.3294545454545454
$ws.Range("G11").Value = 0.1370262390670554
$ws.Range("J11").Value = 0.08454810495626822
$ws.Range("K11").Value = 0.1895043731778426
$ws.Range("L11").Value = 0.5801749271137027
$ws.Range("S11").Value = 0.008746355685131196
$ws.Range("G12").Value = 0.7241379310344828
$ws.Range("J12").Value = 0.2167487684729064
$ws.Range("L12").Value = 0.02463054187192118
$ws.Range("S12").Value = 0.03448275862068965
$ws.Range("F15").Value = 0.02602230483271376
$ws.Range("H15").Value = 0.1189591078066914
$ws.Range("I15").Value = 0.03345724907063197
$ws.Range("J15").Value = 0.3011152416356878
$ws.Range("K15").Value = 0.06319702602230483
$ws.Range("M15").Value = 0.01115241635687732
$ws.Range("O15").Value = 0.05204460966542751
$ws.Range("S15").Value = 0.3940520446096654
$ws.Range("F16").Value = 0.0339622641509434
$ws.Range("H16").Value = 0.1849056603773585
$ws.Range("I16").Value = 0.04150943396226415
$ws.Range("J16").Value = 0.4226415094339623
$ws.Range("K16").Value = 0.1132075471698113
$ws.Range("M16").Value = 0.01132075471698113
$ws.Range("O16").Value = 0.04528301886792453
$ws.Range("S16").Value = 0.1471698113207547
$ws.Range("F17").Value = 0.01219512195121951
$ws.Range("H17").Value = 0.1626016260162602
$ws.Range("I17").Value = 0.05284552845528456
$ws.Range("J17").Value = 0.483739837398374
$ws.Range("K17").Value = 0.09959349593495935
$ws.Range("M17").Value = 0.02439024390243903
$ws.Range("O17").Value = 0.05894308943089431
$ws.Range("S17").Value = 0.1056910569105691
$ws.Range("F18").Value = 0.02994011976047904
$ws.Range("H18").Value = 0.1736526946107785
$ws.Range("I18").Value = 0.07784431137724551
$ws.Range("J18").Value = 0.4311377245508982
$ws.Range("K18").Value = 0.1317365269461078
$ws.Range("M18").Value = 0.005988023952095809
$ws.Range("O18").Value = 0.03592814371257485
$ws.Range("S18").Value = 0.1137724550898204
$ws.Range("F19").Value = 0.01885310290652003
$ws.Range("H19").Value = 0.2411626080125687
$ws.Range("I19").Value = 0.04948939512961508
$ws.Range("J19").Value = 0.3770620581304006
$ws.Range("K19").Value = 0.1217596229379419
$ws.Range("M19").Value = 0.01728201099764336
$ws.Range("O19").Value = 0.06284367635506677
$ws.Range("S19").Value = 0.1115475255302435
